# lipidcane_spearman_1_agile.xlsx - "further work on lipidcane2g"
#
# Summary of the edit:
#  - The Spearman-correlation table is recomputed against new/renamed
#    metrics (production units changed from per-ton to per-year totals,
#    and a new "Heat exchanger network error [%]" metric column is added).
#  - The "Fermentation" parameter group (last row of the table) is removed.
#  - The lipidcane parameter group's rows are reordered/relabeled:
#      Capacity -> Lipid content -> Efficiency -> Lipid retention
#    becomes
#      Lipid content -> Lipid retention -> Additional lipid extraction
#      efficiency -> Capacity
#  - All Spearman correlation coefficients are updated to new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the last row of the table (the "Fermentation" / "Solids
#    loading [%]" parameter). This shifts nothing else because it is the
#    last row (row 14) of the sheet.
# ---------------------------------------------------------------------
$ws.Rows.Item(14).Delete()

# ---------------------------------------------------------------------
# 2. Add the new column J ("Heat exchanger network error [%]") and give
#    its header cells (J1, J2) the same formatting as the existing
#    bold/bordered header cells in column I.
# ---------------------------------------------------------------------
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Update row 2 (metric header row) text for the renamed metrics, and
#    set the new column J header text.
# ---------------------------------------------------------------------
$ws.Range("D2").Value = "Biodiesel production [MMGal/yr]"
$ws.Range("E2").Value = "Ethanol production [MMGal/yr]"
$ws.Range("F2").Value = "Electricity production [MMWhr/yr]"
$ws.Range("G2").Value = "Natural gas consumption [MMcf/yr]"
$ws.Range("H2").Value = "Productivity [MMGGE/yr]"
$ws.Range("J2").Value = "Heat exchanger network error [%]"

# ---------------------------------------------------------------------
# 4. Update the lipidcane parameter row labels (B4:B7) to their new
#    order, and write the new Spearman correlation coefficients for
#    every data row (C:J, skipping the always-empty column G), row by
#    row.
# ---------------------------------------------------------------------

# -- Stream-lipidcane group --
$ws.Range("B4").Value = "Lipid content [dry wt. %]"
$ws.Range("C4").Value = 0.05831578188463128
$ws.Range("D4").Value = 0.9696456182418246
$ws.Range("E4").Value = -0.975222476192899
$ws.Range("F4").Value = 0.9008582696983307
$ws.Range("H4").Value = 0.2997154556046182
$ws.Range("I4").Value = 0.7344220184648806
$ws.Range("J4").Value = 0.9747178511291853

$ws.Range("B5").Value = "Lipid retention [%]"
$ws.Range("C5").Value = 0.03243359736134389
$ws.Range("D5").Value = 0.05167930421117217
$ws.Range("E5").Value = 0.002486588163463526
$ws.Range("F5").Value = 0.02176502659860106
$ws.Range("H5").Value = 0.06440701774428069
$ws.Range("I5").Value = 0.02732779376511175
$ws.Range("J5").Value = 0.02134413969549851

$ws.Range("B6").Value = "Additional lipid extraction efficiency [%]"
$ws.Range("C6").Value = 0.04771422210056887
$ws.Range("D6").Value = 0.07955843799833751
$ws.Range("E6").Value = 0.006257323738292948
$ws.Range("F6").Value = -0.04567332989093319
$ws.Range("H6").Value = 0.01550142657205706
$ws.Range("I6").Value = -0.03496981800679272
$ws.Range("J6").Value = 0.006310359231291893

$ws.Range("B7").Value = "Capacity [ton/hr]"
$ws.Range("C7").Value = 0.06227084946683397
$ws.Range("D7").Value = 0.171234376257375
$ws.Range("E7").Value = 0.1901031738441269
$ws.Range("F7").Value = 0.3424159088326363
$ws.Range("H7").Value = 0.8119513938540557
$ws.Range("I7").Value = 0.6532162222246488
$ws.Range("J7").Value = -0.03414283544528255

# -- Stream-ethanol group (Price [USD/gal]) --
$ws.Range("C8").Value = 0.6944015667200626
$ws.Range("D8").Value = 0.0003540920781636831
$ws.Range("E8").Value = -0.01226857450674298
$ws.Range("F8").Value = -0.005716891140675645
$ws.Range("H8").Value = -0.02644223241768929
$ws.Range("I8").Value = -0.01013140562125622
$ws.Range("J8").Value = 0.001995862832744627

# -- Stream-biodiesel group (Price [USD/gal]) --
$ws.Range("C9").Value = 0.4278411499456459
$ws.Range("D9").Value = -0.00877362227094489
$ws.Range("E9").Value = 0.01120768172830727
$ws.Range("F9").Value = -0.004911329476453179
$ws.Range("H9").Value = 0.006724478860979153
$ws.Range("I9").Value = -0.003638879569555182
$ws.Range("J9").Value = -0.01375664753253893

# -- Stream-natural gas group (Price [USD/cf]) --
$ws.Range("C10").Value = 0.02697139461485578
$ws.Range("D10").Value = 0.003825071385002854
$ws.Range("E10").Value = -0.01243446164937846
$ws.Range("F10").Value = 0.004270110698804427
$ws.Range("H10").Value = -0.008363663662546546
$ws.Range("I10").Value = -0.007364543334581732
$ws.Range("J10").Value = 0.00709263868693979

# -- biorefinery group (Electricity price [USD/kWh]) --
$ws.Range("C11").Value = 0.2113271460530858
$ws.Range("D11").Value = 0.0001760622790424911
$ws.Range("E11").Value = -0.001315681492627259
$ws.Range("F11").Value = -0.001323811348952454
$ws.Range("H11").Value = -0.006028763377150534
$ws.Range("I11").Value = 0.002312393372495735
$ws.Range("J11").Value = -0.003385883176979289

# -- biorefinery group (Operating days [day/yr]) --
$ws.Range("C12").Value = 0.1284625938745037
$ws.Range("D12").Value = 0.1313715268548611
$ws.Range("E12").Value = 0.0815810350072414
$ws.Range("F12").Value = 0.2299450315978012
$ws.Range("H12").Value = 0.4708211171688446
$ws.Range("I12").Value = 0.01965921784236871
$ws.Range("J12").Value = 0.0241950540668351

# -- biorefinery group (IRR [%]) --
$ws.Range("C13").Value = -0.3693392228055689
$ws.Range("D13").Value = 0.01216134144645366
$ws.Range("E13").Value = -0.01501264648850586
$ws.Range("F13").Value = 0.007245784417831376
$ws.Range("H13").Value = -0.00708484953139398
$ws.Range("I13").Value = 0.009378053943122156
$ws.Range("J13").Value = 0.01133482879456183

# ---------------------------------------------------------------------
# 5. Extend the merged header cell from C1:I1 to C1:J1 to cover the new
#    column.
# ---------------------------------------------------------------------
$ws.Range("C1:I1").UnMerge()
$ws.Range("C1:J1").Merge()
